$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 560
$ws.Range("F8").Value = 1202
$ws.Range("F9").Value = 971
$ws.Range("F11").Value = 2223
$ws.Range("F12").Value = 685
$ws.Range("F13").Value = 75
$ws.Range("F15").Value = 738
$ws.Range("F18").Value = 1315
$ws.Range("F22").Value = 11
$ws.Range("F23").Value = 1228
$ws.Range("F24").Value = 294
$ws.Range("F25").Value = 418
$ws.Range("F27").Value = 59
$ws.Range("F30").Value = 11
$ws.Range("F32").Value = 259
$ws.Range("F34").Value = 50
$ws.Range("F41").Value = 202
$ws.Range("F42").Value = 14

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 652
$ws.Range("F10").Value = 365
$ws.Range("F15").Value = 245
$ws.Range("F20").Value = 565
$ws.Range("F23").Value = 420
$ws.Range("F26").Value = 178

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 2130
$ws.Range("F7").Value = 834
$ws.Range("F8").Value = 782
$ws.Range("F11").Value = 771
$ws.Range("F12").Value = 110

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 2130
$ws.Range("F10").Value = 834
$ws.Range("F11").Value = 782
$ws.Range("F14").Value = 560
$ws.Range("F17").Value = 1202
$ws.Range("F18").Value = 971
$ws.Range("F22").Value = 652
$ws.Range("F23").Value = 738
$ws.Range("F25").Value = 365
$ws.Range("F27").Value = 1228
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 418
$ws.Range("F35").Value = 259
$ws.Range("F37").Value = 50
$ws.Range("F42").Value = 178
$ws.Range("F49").Value = 202
